$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; everything shifts right by one.
$ws.Columns.Item(1).Insert()

# The newly inserted B1 becomes the "segments" header. Copy the existing
# header formatting (bold/bordered/centered) from a neighboring header cell
# so the style matches exactly, then set the text afterwards.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "segments"

# Column A now holds a 0-based numeric segment index and takes over the
# bold/bordered style the segment-name column used to have (copy that exact
# formatting from the still-styled name cell in the same row), while column
# B (the segment names) reverts to the plain/default style.
for ($i = 0; $i -lt 19; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).ClearFormats()
}
